$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.645.54'
$ws.Range("E2").Value = '  +1.42%  '
$ws.Range("D3").Value = '''1.631.28'
$ws.Range("E3").Value = '  +1.29%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''213.04'
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("E7").Value = '  +1.18%  '
$ws.Range("E8").Value = '  +1.10%  '
$ws.Range("D9").Value = '''0.0623'
$ws.Range("E9").Value = '  +1.36%  '
$ws.Range("D10").Value = '''18.96'
$ws.Range("E10").Value = '  +2.83%  '
$ws.Range("D11").Value = '''0.0842'
$ws.Range("E11").Value = '  +3.38%  '
$ws.Range("E12").Value = '  +1.40%  '
$ws.Range("D13").Value = '''1.642.69'
$ws.Range("E13").Value = '  +1.82%  '
$ws.Range("D14").Value = '''4.07'
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("D15").Value = '''0.524'
$ws.Range("E15").Value = '  +1.66%  '
$ws.Range("D16").Value = '''26.653.88'
$ws.Range("E16").Value = '  +1.43%  '
$ws.Range("D17").Value = '''63.03'
$ws.Range("E17").Value = '  +1.46%  '
$ws.Range("D18").Value = '''0.0₃0738'
$ws.Range("E18").Value = '  +1.53%  '
$ws.Range("D19").Value = '''209.71'
$ws.Range("E19").Value = '  +4.23%  '
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("E21").Value = '  +0.50%  '
$ws.Range("D22").Value = '''9.43'
$ws.Range("E22").Value = '  +1.12%  '
$ws.Range("D23").Value = '''6.18'
$ws.Range("E23").Value = '  +2.50%  '
$ws.Range("E24").Value = '  +1.25%  '
$ws.Range("D25").Value = '''146.10'
$ws.Range("E25").Value = '  +1.91%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("E27").Value = '  -0.34%  '
$ws.Range("E28").Value = '  +4.75%  '
$ws.Range("E29").Value = '  +1.03%  '
$ws.Range("E30").Value = '  +3.80%  '
$ws.Range("E31").Value = '  -0.35%  '
$ws.Range("E32").Value = '  +1.86%  '
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("E34").Value = '  +1.26%  '
$ws.Range("E35").Value = '  -0.52%  '
$ws.Range("D36").Value = '''1.165.72'
$ws.Range("E36").Value = '  +0.64%  '
$ws.Range("E37").Value = '  +2.17%  '
$ws.Range("E38").Value = '  +2.63%  '
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("D40").Value = '''0.503'
$ws.Range("E40").Value = '  +1.46%  '
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("D42").Value = '''0.791'
$ws.Range("E42").Value = '  +0.90%  '
$ws.Range("D43").Value = '''5.36'
$ws.Range("E43").Value = '  +0.37%  '
$ws.Range("D44").Value = '''1.772.04'
$ws.Range("E44").Value = '  +1.43%  '
$ws.Range("D45").Value = '''92.40'
$ws.Range("E45").Value = '  -0.42%  '
$ws.Range("E46").Value = '  +0.96%  '
$ws.Range("D47").Value = '''54.50'
$ws.Range("E47").Value = '  +1.22%  '
$ws.Range("E48").Value = '  +0.90%  '
$ws.Range("E51").Value = '  -0.02%  '

# Row 49/50: EnergySwap and Mantle swap places with updated price/volume
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = '''7.55'
$ws.Range("E49").Value = '  +4.68%  '
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = '''0.409'
$ws.Range("E50").Value = '  +0.38%  '
